$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = ' Oct 21 2020'
$ws.Cells.Item(2, 2).Value = ' Abu Dhabi'
$ws.Cells.Item(2, 3).Value = 'RCB won by 8 wickets (with 39 balls remaining)'
$ws.Cells.Item(2, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(2, 5).Value = 'Royal Challengers Bangalore'
$ws.Cells.Item(2, 6).Value = 'Pat Cummins '
$ws.Cells.Item(2, 7).Value = '''4'
$ws.Cells.Item(2, 7).Style = "Normal"
$ws.Cells.Item(2, 8).Value = '''17'
$ws.Cells.Item(2, 8).Style = "Normal"
$ws.Cells.Item(2, 9).Value = '''0'
$ws.Cells.Item(2, 9).Style = "Normal"
$ws.Cells.Item(2, 10).Value = '''0'
$ws.Cells.Item(2, 10).Style = "Normal"
$ws.Cells.Item(2, 11).Value = '''23.52'
$ws.Cells.Item(2, 11).Style = "Normal"

$ws.Cells.Item(3, 1).Value = ' Oct 12 2020'
$ws.Cells.Item(3, 2).Value = ' Sharjah'
$ws.Cells.Item(3, 3).Value = 'RCB won by 82 runs'
$ws.Cells.Item(3, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(3, 5).Value = 'Royal Challengers Bangalore'
$ws.Cells.Item(3, 6).Value = 'Pat Cummins '
$ws.Cells.Item(3, 7).Value = '''1'
$ws.Cells.Item(3, 7).Style = "Normal"
$ws.Cells.Item(3, 8).Value = '''3'
$ws.Cells.Item(3, 8).Style = "Normal"
$ws.Cells.Item(3, 9).Value = '''0'
$ws.Cells.Item(3, 9).Style = "Normal"
$ws.Cells.Item(3, 10).Value = '''0'
$ws.Cells.Item(3, 10).Style = "Normal"
$ws.Cells.Item(3, 11).Value = '''33.33'
$ws.Cells.Item(3, 11).Style = "Normal"

$ws.Cells.Item(4, 1).Value = ' Oct 16 2020'
$ws.Cells.Item(4, 2).Value = ' Abu Dhabi'
$ws.Cells.Item(4, 3).Value = 'Mumbai won by 8 wickets (with 19 balls remaining)'
$ws.Cells.Item(4, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(4, 5).Value = 'Mumbai Indians'
$ws.Cells.Item(4, 6).Value = 'Pat Cummins '
$ws.Cells.Item(4, 7).Value = '''53'
$ws.Cells.Item(4, 7).Style = "Normal"
$ws.Cells.Item(4, 8).Value = '''36'
$ws.Cells.Item(4, 8).Style = "Normal"
$ws.Cells.Item(4, 9).Value = '''5'
$ws.Cells.Item(4, 9).Style = "Normal"
$ws.Cells.Item(4, 10).Value = '''2'
$ws.Cells.Item(4, 10).Style = "Normal"
$ws.Cells.Item(4, 11).Value = '''147.22'
$ws.Cells.Item(4, 11).Style = "Normal"

$ws.Cells.Item(5, 1).Value = ' Oct 3 2020'
$ws.Cells.Item(5, 2).Value = ' Sharjah'
$ws.Cells.Item(5, 3).Value = 'Capitals won by 18 runs'
$ws.Cells.Item(5, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(5, 5).Value = 'Delhi Capitals'
$ws.Cells.Item(5, 6).Value = 'Pat Cummins '
$ws.Cells.Item(5, 7).Value = '''5'
$ws.Cells.Item(5, 7).Style = "Normal"
$ws.Cells.Item(5, 8).Value = '''4'
$ws.Cells.Item(5, 8).Style = "Normal"
$ws.Cells.Item(5, 9).Value = '''1'
$ws.Cells.Item(5, 9).Style = "Normal"
$ws.Cells.Item(5, 10).Value = '''0'
$ws.Cells.Item(5, 10).Style = "Normal"
$ws.Cells.Item(5, 11).Value = '''125.00'
$ws.Cells.Item(5, 11).Style = "Normal"

$ws.Cells.Item(6, 1).Value = ' Oct 24 2020'
$ws.Cells.Item(6, 2).Value = ' Abu Dhabi'
$ws.Cells.Item(6, 3).Value = 'KKR won by 59 runs'
$ws.Cells.Item(6, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(6, 5).Value = 'Delhi Capitals'
$ws.Cells.Item(6, 6).Value = 'Pat Cummins '
$ws.Cells.Item(6, 7).Value = '''0'
$ws.Cells.Item(6, 7).Style = "Normal"
$ws.Cells.Item(6, 8).Value = '''0'
$ws.Cells.Item(6, 8).Style = "Normal"
$ws.Cells.Item(6, 9).Value = '''0'
$ws.Cells.Item(6, 9).Style = "Normal"
$ws.Cells.Item(6, 10).Value = '''0'
$ws.Cells.Item(6, 10).Style = "Normal"
$ws.Cells.Item(6, 11).Value = '''-'
$ws.Cells.Item(6, 11).Style = "Normal"

$ws.Cells.Item(7, 1).Value = ' Oct 7 2020'
$ws.Cells.Item(7, 2).Value = ' Abu Dhabi'
$ws.Cells.Item(7, 3).Value = 'KKR won by 10 runs'
$ws.Cells.Item(7, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(7, 5).Value = 'Chennai Super Kings'
$ws.Cells.Item(7, 6).Value = 'Pat Cummins '
$ws.Cells.Item(7, 7).Value = '''17'
$ws.Cells.Item(7, 7).Style = "Normal"
$ws.Cells.Item(7, 8).Value = '''9'
$ws.Cells.Item(7, 8).Style = "Normal"
$ws.Cells.Item(7, 9).Value = '''1'
$ws.Cells.Item(7, 9).Style = "Normal"
$ws.Cells.Item(7, 10).Value = '''1'
$ws.Cells.Item(7, 10).Style = "Normal"
$ws.Cells.Item(7, 11).Value = '''188.88'
$ws.Cells.Item(7, 11).Style = "Normal"

$ws.Cells.Item(8, 1).Value = ' Oct 26 2020'
$ws.Cells.Item(8, 2).Value = ' Sharjah'
$ws.Cells.Item(8, 3).Value = 'Kings XI won by 8 wickets (with 7 balls remaining)'
$ws.Cells.Item(8, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(8, 5).Value = 'Kings XI Punjab'
$ws.Cells.Item(8, 6).Value = 'Pat Cummins '
$ws.Cells.Item(8, 7).Value = '''1'
$ws.Cells.Item(8, 7).Style = "Normal"
$ws.Cells.Item(8, 8).Value = '''8'
$ws.Cells.Item(8, 8).Style = "Normal"
$ws.Cells.Item(8, 9).Value = '''0'
$ws.Cells.Item(8, 9).Style = "Normal"
$ws.Cells.Item(8, 10).Value = '''0'
$ws.Cells.Item(8, 10).Style = "Normal"
$ws.Cells.Item(8, 11).Value = '''12.50'
$ws.Cells.Item(8, 11).Style = "Normal"

$ws.Cells.Item(9, 1).Value = ' Oct 10 2020'
$ws.Cells.Item(9, 2).Value = ' Abu Dhabi'
$ws.Cells.Item(9, 3).Value = 'KKR won by 2 runs'
$ws.Cells.Item(9, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(9, 5).Value = 'Kings XI Punjab'
$ws.Cells.Item(9, 6).Value = 'Pat Cummins '
$ws.Cells.Item(9, 7).Value = '''5'
$ws.Cells.Item(9, 7).Style = "Normal"
$ws.Cells.Item(9, 8).Value = '''4'
$ws.Cells.Item(9, 8).Style = "Normal"
$ws.Cells.Item(9, 9).Value = '''0'
$ws.Cells.Item(9, 9).Style = "Normal"
$ws.Cells.Item(9, 10).Value = '''0'
$ws.Cells.Item(9, 10).Style = "Normal"
$ws.Cells.Item(9, 11).Value = '''125.00'
$ws.Cells.Item(9, 11).Style = "Normal"

$ws.Cells.Item(10, 1).Value = ' Nov 1 2020'
$ws.Cells.Item(10, 2).Value = ' Dubai (DSC)'
$ws.Cells.Item(10, 3).Value = 'KKR won by 60 runs'
$ws.Cells.Item(10, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(10, 5).Value = 'Rajasthan Royals'
$ws.Cells.Item(10, 6).Value = 'Pat Cummins '
$ws.Cells.Item(10, 7).Value = '''15'
$ws.Cells.Item(10, 7).Style = "Normal"
$ws.Cells.Item(10, 8).Value = '''11'
$ws.Cells.Item(10, 8).Style = "Normal"
$ws.Cells.Item(10, 9).Value = '''0'
$ws.Cells.Item(10, 9).Style = "Normal"
$ws.Cells.Item(10, 10).Value = '''1'
$ws.Cells.Item(10, 10).Style = "Normal"
$ws.Cells.Item(10, 11).Value = '''136.36'
$ws.Cells.Item(10, 11).Style = "Normal"

$ws.Cells.Item(11, 1).Value = ' Sep 30 2020'
$ws.Cells.Item(11, 2).Value = ' Dubai (DSC)'
$ws.Cells.Item(11, 3).Value = 'KKR won by 37 runs'
$ws.Cells.Item(11, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(11, 5).Value = 'Rajasthan Royals'
$ws.Cells.Item(11, 6).Value = 'Pat Cummins '
$ws.Cells.Item(11, 7).Value = '''12'
$ws.Cells.Item(11, 7).Style = "Normal"
$ws.Cells.Item(11, 8).Value = '''10'
$ws.Cells.Item(11, 8).Style = "Normal"
$ws.Cells.Item(11, 9).Value = '''1'
$ws.Cells.Item(11, 9).Style = "Normal"
$ws.Cells.Item(11, 10).Value = '''0'
$ws.Cells.Item(11, 10).Style = "Normal"
$ws.Cells.Item(11, 11).Value = '''120.00'
$ws.Cells.Item(11, 11).Style = "Normal"

$ws.Cells.Item(12, 1).Value = ' Sep 23 2020'
$ws.Cells.Item(12, 2).Value = ' Abu Dhabi'
$ws.Cells.Item(12, 3).Value = 'Mumbai won by 49 runs'
$ws.Cells.Item(12, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(12, 5).Value = 'Mumbai Indians'
$ws.Cells.Item(12, 6).Value = 'Pat Cummins '
$ws.Cells.Item(12, 7).Value = '''33'
$ws.Cells.Item(12, 7).Style = "Normal"
$ws.Cells.Item(12, 8).Value = '''12'
$ws.Cells.Item(12, 8).Style = "Normal"
$ws.Cells.Item(12, 9).Value = '''1'
$ws.Cells.Item(12, 9).Style = "Normal"
$ws.Cells.Item(12, 10).Value = '''4'
$ws.Cells.Item(12, 10).Style = "Normal"
$ws.Cells.Item(12, 11).Value = '''275.00'
$ws.Cells.Item(12, 11).Style = "Normal"
